$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "enclos arctique"

$ws.Range("B24").Select()
